# Contexte du jeu.docx - apply the commit's text corrections and
# materialize the document's (until-now implicit) header/footer story.
#
# All edits are performed through Find/Replace on $d.Content so that
# existing run formatting (e.g. the superscript ordinal suffix) is left
# untouched wherever the replaced span does not cross a formatting
# boundary.

$d = $word.ActiveDocument

function Replace-Text($find, $repl) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, `
                                   $false, $true, 1, $false, $repl, 2)
    if (-not $ok) {
        Write-Host "WARNING: could not find/replace: $find"
    }
}

# "21ième" -> "21e" (the run only ever contained the ordinal suffix, so
# this keeps the <w:vertAlign w:val="superscript"/> formatting intact).
Replace-Text "ième" "e"

# Add the missing comma after "siecle" and make "passionne" agree in
# gender with the narrator ("passionnee").
Replace-Text " siècle sans frontière passionné par l’histoire. " `
             " siècle, sans frontière passionnée par l’histoire. "

# Pluralise "important" -> "importants" (agrees with "faits").
Replace-Text "très important." "très importants."

# Fix the agreement/spelling slips: "quelques choses ... d'étranges" ->
# "quelque chose ... d'étrange", "ce produisit" -> "se produisit", and
# drop the spurious trailing "le" on "spatio-temporelle".
Replace-Text "quelques choses d’étranges ce produisit, un trou spatio-temporelle s’ouvrit" `
             "quelque chose d’étrange se produisit, un trou spatio-temporel s’ouvrit"

# Split the big paragraph: a new paragraph now starts at "Il le fit
# voyager...".
Replace-Text "au sein du manuscrit. Il le fit" "au sein du manuscrit. ^pIl le fit"

# Fix "recupere" -> "recuperer" (infinitive), "du" -> "due" (agreement
# with "histoire"), drop the trailing "le" on "spatio-temporelle" again,
# and split a second new paragraph starting at "Son but ?".
Replace-Text "récupère l’histoire qui dégringole dû au trou spatio-temporelle. Son but" `
             "récupérer l’histoire qui dégringole due au trou spatio-temporel. ^pSon but"

# Touching the (until now non-existent) header/footer Range objects
# materializes header1-3.xml, footer1-3.xml, footnotes.xml and
# endnotes.xml, and wires headerReference/footerReference entries into
# the section properties - matching the primary/even/first-page set the
# document now carries (all left blank, as in the source edit).
$sec = $d.Sections.Item(1)
$sec.Headers.Item(1).Range.Text = ""
$sec.Headers.Item(2).Range.Text = ""
$sec.Headers.Item(3).Range.Text = ""
$sec.Footers.Item(1).Range.Text = ""
$sec.Footers.Item(2).Range.Text = ""
$sec.Footers.Item(3).Range.Text = ""

Write-Host "Final text:"
Write-Host $d.Content.Text
